$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 46035
$ws.Range("B2").Value = 12929.7863394535
$ws.Range("C2").Value = 12374.8314682575
$ws.Range("D2").Value = 20355.86
$ws.Range("E2").Value = 8761.64730685041
$ws.Range("F2").Value = 32.525782296162

$ws.Range("A3").Value = 46036
$ws.Range("B3").Value = 12426.2882452703
$ws.Range("C3").Value = 11823.039006589
$ws.Range("D3").Value = 12075.86
$ws.Range("E3").Value = 8451.27092424071
$ws.Range("F3").Value = 341.602080451238

$ws.Range("A4").Value = 46037
$ws.Range("B4").Value = 12328.2953422514
$ws.Range("C4").Value = 11682.7097051479
$ws.Range("D4").Value = 12075.86
$ws.Range("E4").Value = 8385.97774325399
$ws.Range("F4").Value = 333.034477016744

$ws.Range("A5").Value = 46038
$ws.Range("B5").Value = 12125.4136990334
$ws.Range("C5").Value = 10749.455105493
$ws.Range("D5").Value = 12075.86
$ws.Range("E5").Value = 8243.00462102261
$ws.Range("F5").Value = 288.191655271485

$ws.Range("A6").Value = 46039
$ws.Range("B6").Value = 4569.32536066762
$ws.Range("C6").Value = 7163.47319905367
$ws.Range("D6").Value = 12075.86
$ws.Range("E6").Value = 7641.89200886062
$ws.Range("F6").Value = 113.729383663095

$ws.Range("A7").Value = 46040
$ws.Range("B7").Value = 5142.89193121302
$ws.Range("C7").Value = 7254.14751694543
$ws.Range("D7").Value = 12075.86
$ws.Range("E7").Value = 8468.36790424682
$ws.Range("F7").Value = 151.94397588301

$ws.Range("A8").Value = 46041
$ws.Range("B8").Value = 12733.0845017406
$ws.Range("C8").Value = 11027.1593785906
$ws.Range("D8").Value = 12075.86
$ws.Range("E8").Value = 8656.68492000708
$ws.Range("F8").Value = 316.999345774902

$ws.Range("A9").Value = 46042
$ws.Range("B9").Value = 12733.0845017406
$ws.Range("C9").Value = 11654.2686912947
$ws.Range("D9").Value = 12075.86
$ws.Range("E9").Value = 8656.68492000708
$ws.Range("F9").Value = 343.128900470908

$ws.Range("A10").Value = 46043
$ws.Range("B10").Value = 12733.0845017406
$ws.Range("C10").Value = 11883.3228163454
$ws.Range("D10").Value = 12075.86
$ws.Range("E10").Value = 8656.68492000708
$ws.Range("F10").Value = 352.672822348021

$ws.Range("A11").Value = 46044
$ws.Range("B11").Value = 12733.0845017406
$ws.Range("C11").Value = 11891.2763992722
$ws.Range("D11").Value = 12075.86
$ws.Range("E11").Value = 8656.68492000708
$ws.Range("F11").Value = 353.004221636637

$ws.Range("A12").Value = 46045
$ws.Range("B12").Value = 12733.0845017406
$ws.Range("C12").Value = 11290.8080693994
$ws.Range("D12").Value = 12075.86
$ws.Range("E12").Value = 8656.68492000708
$ws.Range("F12").Value = 327.984707891937

$ws.Range("A13").Value = 46046
$ws.Range("B13").Value = 5107.54318705847
$ws.Range("C13").Value = 8159.00460886081
$ws.Range("D13").Value = 12075.86
$ws.Range("E13").Value = 8263.53896066229
$ws.Range("F13").Value = 181.111815396796

$ws.Range("A14").Value = 46047
$ws.Range("B14").Value = 5000.01932310789
$ws.Range("C14").Value = 8311.70165033068
$ws.Range("D14").Value = 12075.86
$ws.Range("E14").Value = 8255.79441006063
$ws.Range("F14").Value = 187.151502516305

$ws.Range("A15").Value = 46048
$ws.Range("B15").Value = 12315.9682835607
$ws.Range("C15").Value = 12180.5124607031
$ws.Range("D15").Value = 12075.86
$ws.Range("E15").Value = 8313.37238908991
$ws.Range("F15").Value = 350.751035408043

